$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (2023-10-05 -> 2023-10-06, serial 45204 -> 45205) for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$ws.Range("C2:C$lastRow").Value = 45205
